$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The table gained a new "2022" column (P), extending the existing
# 2010-2021 year columns (D:O). Bring column P's formatting in line with
# column O (the previous last column) first, then fill in the new values.
$ws.Range("O3:O14").Copy() | Out-Null
$ws.Range("P3:P14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("P4").Value = 2022
$ws.Range("P5").Value = 96.969944810665083
$ws.Range("P6").Value = 96.173557859042035
$ws.Range("P7").Value = 62.289845326160055
$ws.Range("P8").Value = 100
$ws.Range("P9").Value = 100
$ws.Range("P10").Value = "-"
$ws.Range("P11").Value = 100
$ws.Range("P12").Value = 58.090784503861151
$ws.Range("P13").Value = 100
$ws.Range("P14").Value = 100

# Matches the dimension/selection recorded after the edit.
$ws.Range("Q4").Select() | Out-Null
